# Commit: "Better UI and upfront reports"
# - Refresh the "Generated:" timestamp on the Overview sheet.
# - Add three new report sheets at the end of the workbook:
#     "Member Expense Breakdown", "Upfront Payments", "Payment Summary"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while preserving it as literal text.
# Excel's Range.Value setter auto-coerces strings that look like plain
# numbers or bare dates (e.g. "13.00", "9/15/2025") into real numeric /
# date values. The source workbook stores every one of these report
# values as literal text (t="str"), so any cell whose text would
# otherwise be re-interpreted gets its NumberFormat pinned to "@" (Text)
# first, which makes Excel keep the exact original string.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Cell, [string]$Value)

    if ($Value -match '^-?\d+(\.\d+)?$' -or $Value -match '^\d{1,2}/\d{1,2}/\d{4}$') {
        $Cell.NumberFormat = "@"
    }
    $Cell.Value = $Value
}

# ---------------------------------------------------------------------------
# Helper: fill a whole sheet from a tab-separated block of text (first line
# is the header row). Every cell is written as text, matching the source
# file's t="str" cells.
# ---------------------------------------------------------------------------
function Fill-SheetFromTsv {
    param($Sheet, [string]$Tsv)

    $rows = $Tsv -split "`n"
    for ($r = 0; $r -lt $rows.Length; $r++) {
        $cols = $rows[$r] -split "`t"
        for ($c = 0; $c -lt $cols.Length; $c++) {
            $cell = $Sheet.Cells.Item($r + 1, $c + 1)
            Set-TextValue $cell $cols[$c]
        }
    }
}

# ---------------------------------------------------------------------------
# 1. Bump the "Generated:" timestamp on the Overview sheet (B2).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "9/15/2025 9:00:29 PM"

# ---------------------------------------------------------------------------
# 2. Add "Member Expense Breakdown" (sheet8) after "Gig Profitability".
# ---------------------------------------------------------------------------
$memberExpenseBreakdownTsv = @"
Member`tFood Individual (€)`tTravel Equal (€)`tTravel Group (€)`tAdvertisement (€)`tOther (€)`tPaid Out Total (€)`tNet Expense Share (€)`tIncome (€)`tNet Balance (€)`tNet + BTW (€)
Max`t13.00`t36.42`t0.00`t7.89`t0.00`t140.30`t57.31`t152.94`t236.72`t258.02
Naut`t13.00`t36.42`t0.00`t7.89`t0.00`t0.00`t57.31`t152.94`t96.42`t105.10
Filip`t0.00`t36.42`t0.00`t7.89`t0.00`t60.00`t44.31`t76.47`t92.95`t101.31
Dani`t13.00`t36.42`t0.00`t7.89`t0.00`t71.00`t57.31`t152.94`t167.42`t182.49
Pedro`t0.00`t36.42`t0.00`t7.89`t0.00`t45.70`t44.31`t152.94`t155.12`t169.08
Roman`t0.00`t36.42`t0.00`t7.89`t0.00`t9.62`t44.31`t76.47`t42.57`t46.40
Frans`t13.00`t36.42`t0.00`t7.89`t0.00`t0.00`t57.31`t76.47`t19.95`t21.74
Steve`t0.00`t36.42`t0.00`t7.89`t0.00`t42.60`t44.31`t152.94`t152.02`t165.70
Kimo`t0.00`t36.42`t0.00`t7.89`t0.00`t66.00`t44.31`t152.94`t175.42`t191.21
Bandpot`t0.00`t36.42`t0.00`t0.00`t0.00`t0.00`t36.42`t152.94`t109.42`t119.27
"@.Trim()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$memberExpenseBreakdown = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$memberExpenseBreakdown.Name = "Member Expense Breakdown"
Fill-SheetFromTsv $memberExpenseBreakdown $memberExpenseBreakdownTsv

# ---------------------------------------------------------------------------
# 3. Add "Upfront Payments" (sheet9) after "Member Expense Breakdown".
# ---------------------------------------------------------------------------
$upfrontPaymentsTsv = @"
Member`tExpense Type`tDescription`tAmount Paid (€)`tAssociated Gigs`tDate Added`tSplit Among`tMember Share (€)
Steve`tcar-travel`tCar travel for Delft: 121km @ €0.20/km`t24.20`tDelft`t9/15/2025`t10 members (equal split)`t2.42
Kimo`tcar-travel`tCar travel for Delft: 130km @ €0.20/km`t26.00`tDelft`t9/15/2025`t10 members (equal split)`t2.60
Max`tcar-travel`tCar travel for Delft: 140km @ €0.20/km`t28.00`tDelft`t9/15/2025`t10 members (equal split)`t2.80
Pedro`tgeneral-travel`tTravel for Delft: Train`t31.40`tDelft`t9/15/2025`t10 members (equal split)`t3.14
Max`tcar-travel`tCar travel for Camping: 230km @ €0.20/km`t46.00`tCamping`t9/15/2025`t10 members (equal split)`t4.60
Steve`tcar-travel`tCar travel for Camping: 92km @ €0.20/km`t18.40`tCamping`t9/15/2025`t10 members (equal split)`t1.84
Kimo`tcar-travel`tCar travel for Camping: 200km @ €0.20/km`t40.00`tCamping`t9/15/2025`t10 members (equal split)`t4.00
Pedro`tgeneral-travel`tTravel for Camping: Train`t14.30`tCamping`t9/15/2025`t10 members (equal split)`t1.43
Roman`tgeneral-travel`tTravel for Camping: Train`t9.62`tCamping`t9/15/2025`t10 members (equal split)`t0.96
Max`tparking`tParking for Camping: Parking camping`t9.50`tCamping`t9/15/2025`t10 members (equal split)`t0.95
Max`tcar-travel`tCar travel for Camping: Fixed amount`t56.80`tCamping`t9/15/2025`t10 members (equal split)`t5.68
Filip`tcar-travel`tCar travel for Camping: Fixed amount`t60.00`tCamping`t9/15/2025`t10 members (equal split)`t6.00
Dani`tadvertisement`tAdvertisement for Delft, Camping: Ads`t71.00`tDelft; Camping`t9/15/2025`tMax, Naut, Filip, Dani, Pedro, Roman, Frans, Steve, Kimo`t7.89
"@.Trim()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$upfrontPayments = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$upfrontPayments.Name = "Upfront Payments"
Fill-SheetFromTsv $upfrontPayments $upfrontPaymentsTsv

# ---------------------------------------------------------------------------
# 4. Add "Payment Summary" (sheet10) after "Upfront Payments".
# ---------------------------------------------------------------------------
$paymentSummaryTsv = @"
Member`tTotal Paid Upfront (€)`tTotal Expense Share (€)`tNet Payment Impact (€)`tIncome (€)`tFinal Balance (€)`tFinal + BTW (€)
Max`t140.30`t56.52`t83.78`t152.94`t236.72`t258.02
Naut`t0.00`t56.52`t-56.52`t152.94`t96.42`t105.10
Filip`t60.00`t43.52`t16.48`t76.47`t92.95`t101.31
Dani`t71.00`t56.52`t14.48`t152.94`t167.42`t182.49
Pedro`t45.70`t43.52`t2.18`t152.94`t155.12`t169.08
Roman`t9.62`t43.52`t-33.90`t76.47`t42.57`t46.40
Frans`t0.00`t56.52`t-56.52`t76.47`t19.95`t21.74
Steve`t42.60`t43.52`t-0.92`t152.94`t152.02`t165.70
Kimo`t66.00`t43.52`t22.48`t152.94`t175.42`t191.21
Bandpot`t0.00`t43.52`t-43.52`t152.94`t109.42`t119.27
"@.Trim()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$paymentSummary = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$paymentSummary.Name = "Payment Summary"
Fill-SheetFromTsv $paymentSummary $paymentSummaryTsv

Write-Output "Sheets now: $($wb.Worksheets.Count)"
